$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Itga11"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.423576
$ws.Range("H2").Value = 19.270728
$ws.Range("I2").Value = 0.001681024218962088
$ws.Range("J2").Value = 0.001681024218962088
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09233266666666667
$ws.Range("N2").Value = 0.276998
$ws.Range("O2").Value = 0.002874858548413657
$ws.Range("P2").Value = 0.002874858548413658
$ws.Range("Q2").Value = 0.5931059016160001
$ws.Range("R2").Value = 5.337953114544
$ws.Range("S2").Value = 0.000004832706845973551
$ws.Range("T2").Value = 0.000004832706845973552

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Itga11"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.423576
$ws.Range("H3").Value = 19.270728
$ws.Range("I3").Value = 0.001681024218962088
$ws.Range("J3").Value = 0.001681024218962088
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 31.51785566666667
$ws.Range("N3").Value = 94.553567
$ws.Range("O3").Value = 0.9813360759751099
$ws.Range("P3").Value = 0.98133607597511
$ws.Range("Q3").Value = 202.457341231864
$ws.Range("R3").Value = 1822.116071086776
$ws.Range("S3").Value = 0.00164964971065538
$ws.Range("T3").Value = 0.00164964971065538

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Itga11"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.423576
$ws.Range("H4").Value = 19.270728
$ws.Range("I4").Value = 0.001681024218962088
$ws.Range("J4").Value = 0.001681024218962088
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02446166666666667
$ws.Range("N4").Value = 0.07338500000000001
$ws.Range("O4").Value = 0.0007616354434881705
$ws.Range("P4").Value = 0.0007616354434881706
$ws.Range("Q4").Value = 0.15713137492
$ws.Range("R4").Value = 1.41418237428
$ws.Range("S4").Value = 0.000001280327626523545
$ws.Range("T4").Value = 0.000001280327626523546

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Itga11"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.423576
$ws.Range("H5").Value = 19.270728
$ws.Range("I5").Value = 0.001681024218962088
$ws.Range("J5").Value = 0.001681024218962088
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4826403333333333
$ws.Range("N5").Value = 1.447921
$ws.Range("O5").Value = 0.01502743003298815
$ws.Range("P5").Value = 0.01502743003298815
$ws.Range("Q5").Value = 3.100276861832
$ws.Range("R5").Value = 27.902491756488
$ws.Range("S5").Value = 0.00002526147383421133
$ws.Range("T5").Value = 0.00002526147383421133

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Itga11"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3580.644531333333
$ws.Range("H6").Value = 10741.933594
$ws.Range("I6").Value = 0.9370403925578976
$ws.Range("J6").Value = 0.9370403925578976
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.09233266666666667
$ws.Range("N6").Value = 0.276998
$ws.Range("O6").Value = 0.002874858548413657
$ws.Range("P6").Value = 0.002874858548413658
$ws.Range("Q6").Value = 330.6104579634236
$ws.Range("R6").Value = 2975.494121670812
$ws.Range("S6").Value = 0.002693858582753961
$ws.Range("T6").Value = 0.002693858582753962

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Itga11"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3580.644531333333
$ws.Range("H7").Value = 10741.933594
$ws.Range("I7").Value = 0.9370403925578976
$ws.Range("J7").Value = 0.9370403925578976
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 31.51785566666667
$ws.Range("N7").Value = 94.553567
$ws.Range("O7").Value = 0.9813360759751099
$ws.Range("P7").Value = 0.98133607597511
$ws.Range("Q7").Value = 112854.2375322033
$ws.Range("R7").Value = 1015688.13778983
$ws.Range("S7").Value = 0.9195515418629437
$ws.Range("T7").Value = 0.9195515418629439

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col1a2"
$ws.Range("C8").Value = "Itga11"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3580.644531333333
$ws.Range("H8").Value = 10741.933594
$ws.Range("I8").Value = 0.9370403925578976
$ws.Range("J8").Value = 0.9370403925578976
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02446166666666667
$ws.Range("N8").Value = 0.07338500000000001
$ws.Range("O8").Value = 0.0007616354434881705
$ws.Range("P8").Value = 0.0007616354434881706
$ws.Range("Q8").Value = 87.5885329772989
$ws.Range("R8").Value = 788.2967967956901
$ws.Range("S8").Value = 0.0007136831749521638
$ws.Range("T8").Value = 0.0007136831749521639

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col1a2"
$ws.Range("C9").Value = "Itga11"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3580.644531333333
$ws.Range("H9").Value = 10741.933594
$ws.Range("I9").Value = 0.9370403925578976
$ws.Range("J9").Value = 0.9370403925578976
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4826403333333333
$ws.Range("N9").Value = 1.447921
$ws.Range("O9").Value = 0.01502743003298815
$ws.Range("P9").Value = 0.01502743003298815
$ws.Range("Q9").Value = 1728.163470150897
$ws.Range("R9").Value = 15553.47123135807
$ws.Range("S9").Value = 0.01408130893724755
$ws.Range("T9").Value = 0.01408130893724756

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Col1a2"
$ws.Range("C10").Value = "Itga11"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9157713333333334
$ws.Range("H10").Value = 2.747314
$ws.Range("I10").Value = 0.0002396537054071653
$ws.Range("J10").Value = 0.0002396537054071653
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.09233266666666667
$ws.Range("N10").Value = 0.276998
$ws.Range("O10").Value = 0.002874858548413657
$ws.Range("P10").Value = 0.002874858548413658
$ws.Range("Q10").Value = 0.08455560926355557
$ws.Range("R10").Value = 0.7610004833720001
$ws.Range("S10").Value = 0.0000006889705036487974
$ws.Range("T10").Value = 0.0000006889705036487976

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col1a2"
$ws.Range("C11").Value = "Itga11"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.9157713333333334
$ws.Range("H11").Value = 2.747314
$ws.Range("I11").Value = 0.0002396537054071653
$ws.Range("J11").Value = 0.0002396537054071653
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 31.51785566666667
$ws.Range("N11").Value = 94.553567
$ws.Range("O11").Value = 0.9813360759751099
$ws.Range("P11").Value = 0.98133607597511
$ws.Range("Q11").Value = 28.86314870767089
$ws.Range("R11").Value = 259.768338369038
$ws.Range("S11").Value = 0.0002351808268571625
$ws.Range("T11").Value = 0.0002351808268571626

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col1a2"
$ws.Range("C12").Value = "Itga11"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9157713333333334
$ws.Range("H12").Value = 2.747314
$ws.Range("I12").Value = 0.0002396537054071653
$ws.Range("J12").Value = 0.0002396537054071653
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02446166666666667
$ws.Range("N12").Value = 0.07338500000000001
$ws.Range("O12").Value = 0.0007616354434881705
$ws.Range("P12").Value = 0.0007616354434881706
$ws.Range("Q12").Value = 0.02240129309888889
$ws.Range("R12").Value = 0.20161163789
$ws.Range("S12").Value = 0.0000001825287562013697
$ws.Range("T12").Value = 0.0000001825287562013697

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col1a2"
$ws.Range("C13").Value = "Itga11"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9157713333333334
$ws.Range("H13").Value = 2.747314
$ws.Range("I13").Value = 0.0002396537054071653
$ws.Range("J13").Value = 0.0002396537054071653
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4826403333333333
$ws.Range("N13").Value = 1.447921
$ws.Range("O13").Value = 0.01502743003298815
$ws.Range("P13").Value = 0.01502743003298815
$ws.Range("Q13").Value = 0.4419881815771111
$ws.Range("R13").Value = 3.977893634194
$ws.Range("S13").Value = 0.000003601379290152529
$ws.Range("T13").Value = 0.00000360137929015253

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col1a2"
$ws.Range("C14").Value = "Itga11"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 233.243637
$ws.Range("H14").Value = 699.7309110000001
$ws.Range("I14").Value = 0.0610389295177331
$ws.Range("J14").Value = 0.06103892951773311
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.09233266666666667
$ws.Range("N14").Value = 0.276998
$ws.Range("O14").Value = 0.002874858548413657
$ws.Range("P14").Value = 0.002874858548413658
$ws.Range("Q14").Value = 21.536006987242
$ws.Range("R14").Value = 193.8240628851781
$ws.Range("S14").Value = 0.0001754782883100737
$ws.Range("T14").Value = 0.0001754782883100738

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col1a2"
$ws.Range("C15").Value = "Itga11"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 233.243637
$ws.Range("H15").Value = 699.7309110000001
$ws.Range("I15").Value = 0.0610389295177331
$ws.Range("J15").Value = 0.06103892951773311
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 31.51785566666667
$ws.Range("N15").Value = 94.553567
$ws.Range("O15").Value = 0.9813360759751099
$ws.Range("P15").Value = 0.98133607597511
$ws.Range("Q15").Value = 7351.339286134394
$ws.Range("R15").Value = 66162.05357520955
$ws.Range("S15").Value = 0.0598997035746535
$ws.Range("T15").Value = 0.05989970357465352

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col1a2"
$ws.Range("C16").Value = "Itga11"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 233.243637
$ws.Range("H16").Value = 699.7309110000001
$ws.Range("I16").Value = 0.0610389295177331
$ws.Range("J16").Value = 0.06103892951773311
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02446166666666667
$ws.Range("N16").Value = 0.07338500000000001
$ws.Range("O16").Value = 0.0007616354434881705
$ws.Range("P16").Value = 0.0007616354434881706
$ws.Range("Q16").Value = 5.705528100415002
$ws.Range("R16").Value = 51.34975290373501
$ws.Range("S16").Value = 0.00004648941215328183
$ws.Range("T16").Value = 0.00004648941215328184

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col1a2"
$ws.Range("C17").Value = "Itga11"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 233.243637
$ws.Range("H17").Value = 699.7309110000001
$ws.Range("I17").Value = 0.0610389295177331
$ws.Range("J17").Value = 0.06103892951773311
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4826403333333333
$ws.Range("N17").Value = 1.447921
$ws.Range("O17").Value = 0.01502743003298815
$ws.Range("P17").Value = 0.01502743003298815
$ws.Range("Q17").Value = 112.572786709559
$ws.Range("R17").Value = 1013.155080386031
$ws.Range("S17").Value = 0.0009172582426162292
$ws.Range("T17").Value = 0.0009172582426162294
